$d = $word.ActiveDocument

# Locate the anchor paragraph via Find on its distinctive text.
$anchor = $d.Content
$anchor.Find.Execute("Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$newLines = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

$insertAfter = $anchor
$insertAfter.Collapse(0)

foreach ($line in $newLines) {
    $insertAfter.InsertParagraphAfter()
    $insertAfter.Collapse(0)
    $insertAfter.MoveStart(1, 1)
    $insertAfter.InsertBefore($line)
    $insertAfter.Collapse(0)
}

Write-Output "done"
